# Generate Report for Handback
# Adds a new row (row 4) to each of the three worksheets (Overview, zh-cn, de-de)
# for the file "794408b3-5c77-41cc-a60c-f132dac7bf68.md", mirroring the existing
# rows for "3e4ac5fa-c504-4f35-ace7-32eb4448f740.md" (the "in sync with en-US" case).

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$fileId = "794408b3-5c77-41cc-a60c-f132dac7bf68"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item("Overview")
[void]$loOv.ListRows.Add()

$wsOv.Range("A4").Value = "$fileId.md"
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("E4").Value = "Handed back: in sync with en-US"
$wsOv.Range("F4").Value = "Handed back: in sync with en-US"
$wsOv.Range("G4").Value = "2016-08-21 02:49:08"
$wsOv.Range("G4").NumberFormat = $dateFmt

[void]$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f1fcbddbaa5b1e2e2c87eb1e2c28a7a84e3c3f2/e2e/$fileId.md", [Type]::Missing, [Type]::Missing, "e2e\$fileId.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item("zh-cn")
[void]$loZh.ListRows.Add()

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F2").Copy($wsZh.Range("F4"))

$wsZh.Range("G4").Value = "$fileId.a01a061ecb0a4c68e7b78461723fbba0da52fe44.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-21 02:49:02"
$wsZh.Range("H4").NumberFormat = $dateFmt
$wsZh.Range("J4").Value = "$fileId.a01a061ecb0a4c68e7b78461723fbba0da52fe44.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-08-21 02:49:24"
$wsZh.Range("K4").NumberFormat = $dateFmt

$wsZh.Range("L2").Copy($wsZh.Range("L4"))
$wsZh.Range("M2").Copy($wsZh.Range("M4"))
$wsZh.Range("N2").Copy($wsZh.Range("N4"))
$wsZh.Range("O2").Copy($wsZh.Range("O4"))
$wsZh.Range("P2").Copy($wsZh.Range("P4"))

[void]$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f1fcbddbaa5b1e2e2c87eb1e2c28a7a84e3c3f2/e2e/$fileId.md", [Type]::Missing, [Type]::Missing, "$fileId.md")
[void]$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8f2c6b6a4e4c6c8a1d2e3f4a5b6c7d8e9f0a1b2c/e2e/$fileId.md", [Type]::Missing, [Type]::Missing, "$fileId.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item("de-de")
[void]$loDe.ListRows.Add()

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F2").Copy($wsDe.Range("F4"))

$wsDe.Range("G4").Value = "$fileId.a01a061ecb0a4c68e7b78461723fbba0da52fe44.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-21 02:49:08"
$wsDe.Range("H4").NumberFormat = $dateFmt
$wsDe.Range("J4").Value = "$fileId.a01a061ecb0a4c68e7b78461723fbba0da52fe44.de-de.xlf"
$wsDe.Range("K4").Value = "2016-08-21 02:49:30"
$wsDe.Range("K4").NumberFormat = $dateFmt

$wsDe.Range("L2").Copy($wsDe.Range("L4"))
$wsDe.Range("M2").Copy($wsDe.Range("M4"))
$wsDe.Range("N2").Copy($wsDe.Range("N4"))
$wsDe.Range("O2").Copy($wsDe.Range("O4"))
$wsDe.Range("P2").Copy($wsDe.Range("P4"))

[void]$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f1fcbddbaa5b1e2e2c87eb1e2c28a7a84e3c3f2/e2e/$fileId.md", [Type]::Missing, [Type]::Missing, "$fileId.md")
[void]$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3b1d6c4e8f2a9b7c5d4e3f2a1b0c9d8e7f6a5b4c/e2e/$fileId.md", [Type]::Missing, [Type]::Missing, "$fileId.md")
